$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 420
$ws.Range("K5").Value = 174
$ws.Range("K6").Value = 246
